# Update cryptos price/volume(1h) data for Thu Feb  2 08:44:24 UTC 2023 GitHub Actions refresh.
# Price (D) and Volume(1h) (E) columns are stored as plain text in the sheet, so we force
# a text number format before writing the value to avoid Excel auto-converting the
# numeric/percentage-looking strings into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "329.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "7.16%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.77%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.261"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.51%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08105"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.88%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.535"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.34%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.637"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.76%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.94%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.34%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9353"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.87%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1326"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "22.17%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1965"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.60%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09108"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.96%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03563"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "7.11%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09576"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.19%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001326"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.18%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006068"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "6.57%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.364"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-5.56%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3516"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.20%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.978"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "10.96%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1342"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.28%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04411"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.53%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001222"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.75%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004301"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.14%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001190"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-8.73%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003990"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.05%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02512"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "15.64%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05179"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.23%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007723"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.67%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1428"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.93%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009212"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.35%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002160"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.80%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01023"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "28.08%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.55%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.28%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003345"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "16.66%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002483"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "147.59%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.28%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.28%"
